$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(234, 44308, 1, 12, 192.2768787053357),
    @(235, 44309, 0, 8, 128.1845858035571),
    @(236, 44310, 2, 9, 144.2076590290018),
    @(237, 44311, 3, 7, 112.1615125781125),
    @(238, 44312, 1, 7, 112.1615125781125)
)

# Copy formatting of the existing A233 cell (date column style) once,
# then paste the formats onto each new date cell below it.
$ws.Range("A233").Copy() | Out-Null

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
